$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix typo in A142: "Checkpoitn 2224" -> "Checkpoint 2224" ---
# (This is the very first new/changed shared string, so it must be set first
#  to land at the correct position in the shared-strings table.)
$ws.Range("A142").Value2 = "Checkpoint 2224"

# --- Copy style (border/number format) from row 148 down to the new rows 149:167 ---
$ws.Range("A148:D148").Copy() | Out-Null
$ws.Range("A149:D167").PasteSpecial(-4122) | Out-Null

# --- Write column-A labels in the exact order the strings were originally authored ---
# (several rows re-use pre-existing strings; the order below reproduces the
#  resulting shared-strings table from the source workbook.)
$ws.Range("A149").Value2 = "Checkpoint 271"
$ws.Range("A150").Value2 = "Checkpoint 596/595"
$ws.Range("A151").Value2 = "Checkoint 872/870"
$ws.Range("A152").Value2 = "Checkpoint 1293/1291"
$ws.Range("A153").Value2 = "Checkpoint 1945"
$ws.Range("A154").Value2 = "Checkpoint 2354"
$ws.Range("A155").Value2 = "Checkpoitn 2941"
$ws.Range("A156").Value2 = "Enter Pipe"
$ws.Range("A157").Value2 = "Get Flag"
$ws.Range("A158").Value2 = "End Level"
$ws.Range("A159").Value2 = "Enter 8-6"
$ws.Range("A160").Value2 = "1st Move"
$ws.Range("A161").Value2 = "Land 1st Orange plat"
$ws.Range("A166").Value2 = "Jump block"
$ws.Range("A162").Value2 = "Jump Pipe"
$ws.Range("A165").Value2 = "Push on spring"
$ws.Range("A164").Value2 = "Checkpoint"
$ws.Range("A167").Value2 = "Spring off 2nd spring (sparks)"

# --- Fill in the timing data (columns B/C) for the new rows ---
$ws.Range("B149").Value2 = 50797
$ws.Range("C149").Value2 = 59960
$ws.Range("B150").Value2 = 50904
$ws.Range("C150").Value2 = 60068
$ws.Range("B151").Value2 = 50995
$ws.Range("C151").Value2 = 60159
$ws.Range("B152").Value2 = 51135
$ws.Range("C152").Value2 = 60299
$ws.Range("B153").Value2 = 51352
$ws.Range("C153").Value2 = 60517
$ws.Range("B154").Value2 = 51487
$ws.Range("C154").Value2 = 60652
$ws.Range("B155").Value2 = 51681
$ws.Range("C155").Value2 = 60846
$ws.Range("B156").Value2 = 52080
$ws.Range("C156").Value2 = 61248
$ws.Range("B157").Value2 = 52250
$ws.Range("C157").Value2 = 61427
$ws.Range("B158").Value2 = 52768
$ws.Range("C158").Value2 = 61945
$ws.Range("B159").Value2 = 53104
$ws.Range("C159").Value2 = 62610
$ws.Range("B160").Value2 = 53332
$ws.Range("C160").Value2 = 62861
$ws.Range("B161").Value2 = 53477
$ws.Range("C161").Value2 = 63021
$ws.Range("B162").Value2 = 53562
$ws.Range("C162").Value2 = 63110
$ws.Range("B163").Value2 = 53643
$ws.Range("C163").Value2 = 63206
$ws.Range("B164").Value2 = 53813
$ws.Range("C164").Value2 = 63376
$ws.Range("B165").Value2 = 53851
$ws.Range("C165").Value2 = 63425
$ws.Range("B166").Value2 = 53952
$ws.Range("C166").Value2 = 63540
$ws.Range("B167").Value2 = 54007
$ws.Range("C167").Value2 = 63597

# --- Fill column D with the elapsed-time formula for the new rows ---
$ws.Range("D149:D167").Formula = "=IF(B149 >  0,C149-B149, 0)"

# --- Update frozen-pane scroll position and active selection to match new data extent ---
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A2").Select() | Out-Null
$win.FreezePanes = $true
$win.ScrollRow = 151
$ws.Range("B168").Select() | Out-Null

Write-Output "done"
